$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows after row 14 (pushes old rows 15-16 down to 18-19)
$ws.Rows("15:17").Insert() | Out-Null

# Row 17 becomes a clone of row 14 (the "(-)" / Atividade header row), then
# gets its Atividade bumped from 2 to 3 and its description restored to
# "Front-End" (row 14 itself will be repurposed below).
$ws.Range("B14:G14").Copy() | Out-Null
$ws.Range("B17:G17").PasteSpecial() | Out-Null
$ws.Range("D17").Value = 3

# Row 14 keeps its Status/Principal/Atividade/Sub-Atividade/Correcao Bugs but
# its description changes to the new "BD" activity.
$ws.Range("G14").Value = "BD"

# Rows 15-16 were blank inserts that inherited row 14's formatting; strip
# that back to the default (unstyled) look before filling them in.
$ws.Range("C15:G16").ClearFormats() | Out-Null

# Two brand new detail rows under the "BD" activity.
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = "Dados do Cartao Para Transacao"
$ws.Range("B15").Interior.ThemeColor = 2

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = "Transacaoes Incertas"
$ws.Range("B16").Interior.ThemeColor = 2

# Old rows 15/16 (now shifted to 18/19) move from Atividade 2 to Atividade 3.
$ws.Range("D18").Value = 3
$ws.Range("D19").Value = 3

# Expand the table / autofilter to include the three new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:G19")) | Out-Null

$ws.Range("G23").Select() | Out-Null
